$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E7").Value = 16.704
$ws.Range("A9").Value = -21.778
$ws.Range("E12").Value = 17.646
$ws.Range("A18").Value = -22.156
$ws.Range("A20").Value = -19.918
$ws.Range("E26").Value = 16.525
$ws.Range("A27").Value = -22.01
$ws.Range("E27").Value = 16.531
$ws.Range("E29").Value = 16.941
$ws.Range("E37").Value = 16.855
$ws.Range("E38").Value = 16.741
$ws.Range("E51").Value = 16.65
$ws.Range("E55").Value = 16.494
$ws.Range("A69").Value = -21.55700000000001
$ws.Range("E69").Value = 17.438
$ws.Range("E70").Value = 17.524
$ws.Range("A76").Value = -20.043
$ws.Range("A82").Value = -22.152
$ws.Range("E83").Value = 16.886
$ws.Range("E102").Value = 16.724
